$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 37505644
$ws.Cells.Item(32, 9).Value = 50002656
$ws.Cells.Item(32, 11).Value = 50002656
$ws.Cells.Item(32, 13).Value = -50002330

$ws.Cells.Item(62, 8).Value = 1965.7059
$ws.Cells.Item(62, 9).Value = 1948.8667
$ws.Cells.Item(62, 11).Value = 1948.8667
$ws.Cells.Item(62, 13).Value = -1324.8667

$ws.Cells.Item(64, 8).Value = 4109.647
$ws.Cells.Item(64, 9).Value = 3576
$ws.Cells.Item(64, 10).Value = 4332
$ws.Cells.Item(64, 11).Value = 3576
$ws.Cells.Item(64, 12).Value = 4332
$ws.Cells.Item(64, 13).Value = -3328
$ws.Cells.Item(64, 14).Value = -4828

$ws.Cells.Item(65, 8).Value = 1965.7059
$ws.Cells.Item(65, 9).Value = 1948.8667
$ws.Cells.Item(65, 11).Value = 9744.333500000001
$ws.Cells.Item(65, 13).Value = -6624.333500000001

$ws.Cells.Item(67, 8).Value = 4109.647
$ws.Cells.Item(67, 9).Value = 3576
$ws.Cells.Item(67, 10).Value = 4332
$ws.Cells.Item(67, 11).Value = 3576
$ws.Cells.Item(67, 12).Value = 4332
$ws.Cells.Item(67, 13).Value = -2718
$ws.Cells.Item(67, 14).Value = -6048

$ws.Cells.Item(103, 8).Value = 702.4
$ws.Cells.Item(103, 9).Value = 298.25
$ws.Cells.Item(103, 11).Value = 894.75
$ws.Cells.Item(103, 13).Value = -308.75

$ws.Cells.Item(106, 8).Value = 62503550
$ws.Cells.Item(106, 9).Value = 71431490
$ws.Cells.Item(106, 11).Value = 71431490
$ws.Cells.Item(106, 13).Value = -71430859

$ws.Cells.Item(113, 8).Value = 7999.75
$ws.Cells.Item(113, 9).Value = 9394.833000000001
$ws.Cells.Item(113, 10).Value = 3814.5
$ws.Cells.Item(113, 11).Value = 9394.833000000001
$ws.Cells.Item(113, 12).Value = 3814.5
$ws.Cells.Item(113, 13).Value = -6140.833000000001
$ws.Cells.Item(113, 14).Value = -10322.5

$ws.Cells.Item(138, 8).Value = 2045.3396
$ws.Cells.Item(138, 9).Value = 986.5517
$ws.Cells.Item(138, 11).Value = 2959.6551
$ws.Cells.Item(138, 13).Value = 2180.3449

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 4099
$ws.Cells.Item(41, 9).Value = 4099
$ws.Cells.Item(41, 11).Value = 4099
$ws.Cells.Item(41, 13).Value = -3685

$ws.Cells.Item(45, 8).Value = 17442.5
$ws.Cells.Item(45, 9).Value = 33474.75
$ws.Cells.Item(45, 10).Value = 1410.25
$ws.Cells.Item(45, 11).Value = 33474.75
$ws.Cells.Item(45, 12).Value = 1410.25
$ws.Cells.Item(45, 13).Value = -33097.75
$ws.Cells.Item(45, 14).Value = -2164.25

$ws.Cells.Item(46, 8).Value = 4616.5
$ws.Cells.Item(46, 10).Value = 4616.5
$ws.Cells.Item(46, 12).Value = 4616.5
$ws.Cells.Item(46, 14).Value = -5254.5

$ws.Cells.Item(61, 8).Value = 1900.129
$ws.Cells.Item(61, 9).Value = 1908.25
$ws.Cells.Item(61, 10).Value = 1824.3334
$ws.Cells.Item(61, 11).Value = 1908.25
$ws.Cells.Item(61, 12).Value = 1824.3334
$ws.Cells.Item(61, 13).Value = -1696.25
$ws.Cells.Item(61, 14).Value = -2248.3334

$ws.Cells.Item(97, 8).Value = 21741670
$ws.Cells.Item(97, 9).Value = 29414148
$ws.Cells.Item(97, 11).Value = 29414148
$ws.Cells.Item(97, 13).Value = -29413652

$ws.Cells.Item(132, 8).Value = 1516.58
$ws.Cells.Item(132, 9).Value = 1371.7805
$ws.Cells.Item(132, 11).Value = 4115.3415
$ws.Cells.Item(132, 13).Value = -1585.3415

$ws.Cells.Item(136, 8).Value = 1900.129
$ws.Cells.Item(136, 9).Value = 1908.25
$ws.Cells.Item(136, 10).Value = 1824.3334
$ws.Cells.Item(136, 11).Value = 5724.75
$ws.Cells.Item(136, 12).Value = 5473.0002
$ws.Cells.Item(136, 13).Value = -3174.75
$ws.Cells.Item(136, 14).Value = -10573.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(97, 8).Value = 3185
$ws.Cells.Item(97, 9).Value = 3185
$ws.Cells.Item(97, 11).Value = 3185
$ws.Cells.Item(97, 13).Value = -2194

$ws.Cells.Item(105, 8).Value = 3417.2632
$ws.Cells.Item(105, 9).Value = 2297.484
$ws.Cells.Item(105, 11).Value = 2297.484
$ws.Cells.Item(105, 13).Value = -550.4839999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 87.09999999999999
$ws.Cells.Item(7, 9).Value = 65.5
$ws.Cells.Item(7, 11).Value = 65.5
$ws.Cells.Item(7, 13).Value = 47.5

$ws.Cells.Item(10, 8).Value = 3128.75
$ws.Cells.Item(10, 10).Value = 10008
$ws.Cells.Item(10, 12).Value = 10008
$ws.Cells.Item(10, 14).Value = -10286

$ws.Cells.Item(62, 8).Value = 45457492
$ws.Cells.Item(62, 9).Value = 142860690
$ws.Cells.Item(62, 10).Value = 2664.7334
$ws.Cells.Item(62, 11).Value = 142860690
$ws.Cells.Item(62, 12).Value = 2664.7334
$ws.Cells.Item(62, 13).Value = -142860066
$ws.Cells.Item(62, 14).Value = -3912.7334

$ws.Cells.Item(65, 8).Value = 45457492
$ws.Cells.Item(65, 9).Value = 142860690
$ws.Cells.Item(65, 10).Value = 2664.7334
$ws.Cells.Item(65, 11).Value = 714303450
$ws.Cells.Item(65, 12).Value = 13323.667
$ws.Cells.Item(65, 13).Value = -714300330
$ws.Cells.Item(65, 14).Value = -19563.667

$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 1327.92
$ws.Cells.Item(107, 9).Value = 593.25
$ws.Cells.Item(107, 10).Value = 4266.6
$ws.Cells.Item(107, 11).Value = 593.25
$ws.Cells.Item(107, 12).Value = 4266.6
$ws.Cells.Item(107, 13).Value = 1326.75
$ws.Cells.Item(107, 14).Value = -8106.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 907.3158
$ws.Cells.Item(5, 9).Value = 908.75
$ws.Cells.Item(5, 11).Value = 2726.25
$ws.Cells.Item(5, 13).Value = -2614.25

$ws.Cells.Item(7, 8).Value = 113.71429
$ws.Cells.Item(7, 10).Value = 78.5
$ws.Cells.Item(7, 12).Value = 235.5
$ws.Cells.Item(7, 14).Value = -459.5

$ws.Cells.Item(34, 8).Value = 316.5
$ws.Cells.Item(34, 10).Value = 750
$ws.Cells.Item(34, 12).Value = 2250
$ws.Cells.Item(34, 14).Value = -2418

$ws.Cells.Item(39, 8).Value = 8760.875
$ws.Cells.Item(39, 10).Value = 10066.167
$ws.Cells.Item(39, 12).Value = 30198.501
$ws.Cells.Item(39, 14).Value = -30786.501

$ws.Cells.Item(55, 8).Value = 1179.625
$ws.Cells.Item(55, 9).Value = 572.8333
$ws.Cells.Item(55, 10).Value = 3000
$ws.Cells.Item(55, 11).Value = 1718.4999
$ws.Cells.Item(55, 12).Value = 9000
$ws.Cells.Item(55, 13).Value = -1541.4999
$ws.Cells.Item(55, 14).Value = -9354

$ws.Cells.Item(134, 8).Value = 1038.6666
$ws.Cells.Item(134, 9).Value = 1038.6666
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 3115.9998
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = 1954.0002
$ws.Cells.Item(134, 14).ClearContents()

$ws.Cells.Item(135, 8).Value = 907.3158
$ws.Cells.Item(135, 9).Value = 908.75
$ws.Cells.Item(135, 11).Value = 8178.75
$ws.Cells.Item(135, 13).Value = -5643.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 16834.334
$ws.Cells.Item(9, 10).Value = 40008
$ws.Cells.Item(9, 12).Value = 40008
$ws.Cells.Item(9, 14).Value = -40348

$ws.Cells.Item(94, 8).Value = 36497.25
$ws.Cells.Item(94, 10).Value = 36497.25
$ws.Cells.Item(94, 12).Value = 36497.25
$ws.Cells.Item(94, 14).Value = -37849.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 44997.5
$ws.Cells.Item(14, 9).Value = 45000
$ws.Cells.Item(14, 10).Value = 44995
$ws.Cells.Item(14, 11).Value = 45000
$ws.Cells.Item(14, 12).Value = 44995
$ws.Cells.Item(14, 13).Value = -44828
$ws.Cells.Item(14, 14).Value = -45339

$ws.Cells.Item(82, 8).Value = 466.10388
$ws.Cells.Item(82, 9).Value = 466.10388
$ws.Cells.Item(82, 11).Value = 466.10388
$ws.Cells.Item(82, 13).Value = -105.10388

$ws.Cells.Item(85, 8).Value = 466.10388
$ws.Cells.Item(85, 9).Value = 466.10388
$ws.Cells.Item(85, 11).Value = 466.10388
$ws.Cells.Item(85, 13).Value = 781.89612

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 399
$ws.Cells.Item(30, 10).Value = 399
$ws.Cells.Item(30, 12).Value = 399
$ws.Cells.Item(30, 14).Value = -613

$ws.Cells.Item(62, 8).Value = 125008460
$ws.Cells.Item(62, 9).Value = 200007340
$ws.Cells.Item(62, 10).Value = 10333
$ws.Cells.Item(62, 11).Value = 200007340
$ws.Cells.Item(62, 12).Value = 10333
$ws.Cells.Item(62, 13).Value = -200006716
$ws.Cells.Item(62, 14).Value = -11581

$ws.Cells.Item(65, 8).Value = 125008460
$ws.Cells.Item(65, 9).Value = 200007340
$ws.Cells.Item(65, 10).Value = 10333
$ws.Cells.Item(65, 11).Value = 1000036700
$ws.Cells.Item(65, 12).Value = 51665
$ws.Cells.Item(65, 13).Value = -1000033580
$ws.Cells.Item(65, 14).Value = -57905

$ws.Cells.Item(81, 8).Value = 8449502
$ws.Cells.Item(81, 9).Value = 13131361
$ws.Cells.Item(81, 10).Value = 256250
$ws.Cells.Item(81, 11).Value = 26262722
$ws.Cells.Item(81, 12).Value = 512500
$ws.Cells.Item(81, 13).Value = -26261661
$ws.Cells.Item(81, 14).Value = -514622

$ws.Cells.Item(84, 8).Value = 8449502
$ws.Cells.Item(84, 9).Value = 13131361
$ws.Cells.Item(84, 10).Value = 256250
$ws.Cells.Item(84, 11).Value = 131313610
$ws.Cells.Item(84, 12).Value = 2562500
$ws.Cells.Item(84, 13).Value = -131308306
$ws.Cells.Item(84, 14).Value = -2573108

$ws.Cells.Item(136, 8).Value = 30313.555
$ws.Cells.Item(136, 9).Value = 27633.154
$ws.Cells.Item(136, 11).Value = 82899.462
$ws.Cells.Item(136, 13).Value = -80349.462
